$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Insert a new column before H (date) to hold the new "property_category" data.
# This shifts existing H (date) -> I, I (legislator_name) -> J, J (legislator_id) -> K,
# while correctly carrying over cell styles (border/bold for header row, plain for data rows).
$ws.Columns("H").Insert(-4121)

$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H20").Value = "stock"

# Fix stray mid-string spaces in company names (shared strings clean-up)
$ws.Range("B2").Value  = "聯華實業股份有限公司"
$ws.Range("B3").Value  = "華新科技股份有限公司"
$ws.Range("B4").Value  = "佳能企業股份有限公司"
$ws.Range("B5").Value  = "彰化商業銀行股份有限公司"
$ws.Range("B6").Value  = "臺灣中小企業銀行股份有限公司’"
$ws.Range("B7").Value  = "中華開發金融控股股份有限公司"
$ws.Range("B8").Value  = "玉山金融控股股份有限公司"
$ws.Range("B9").Value  = "台新金融控股股份有限公司"
$ws.Range("B10").Value = "第一金融控股股份有限公司"
$ws.Range("B11").Value = "遠東百貨股份有限公司"
$ws.Range("B12").Value = "_創資通股份有限公司"
$ws.Range("B13").Value = "遠雄建設事業股份有限公司"
$ws.Range("B14").Value = "遠雄自賀港投資控股股份有限公司"
$ws.Range("B15").Value = "瑞儀光電股份有限公司"
$ws.Range("B16").Value = "康舒科技股份有限公司"
$ws.Range("B17").Value = "台灣水泥股份有限公司"
$ws.Range("B18").Value = "台灣積體電路製造股份有限公司"
$ws.Range("B19").Value = "華南金融控股股份有限公司"
$ws.Range("B20").Value = "寶來曼氏期貨股份有限公司"
